# Update the "HotStock_Top20" style ranking lists in columns A, B and C.
# Rows 2-21 hold the ranked stock names for each of the three sources
# (财联社 / 东方财富 / 同花顺). The lists have been reshuffled; write the
# new ranking for each column in one shot using the Range.Value array form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(
    "阳光电源",
    "巨轮智能",
    "岩山科技",
    "卧龙电驱",
    "长城军工",
    "中际旭创",
    "利欧股份",
    "国轩高科",
    "吉视传媒",
    "浙江荣泰",
    "华胜天成",
    "山子高科",
    "新易盛",
    "北方稀土",
    "东方财富",
    "寒武纪-U",
    "天普股份",
    "通富微电",
    "春兴精工",
    "至纯科技"
)

$colB = @(
    "岩山科技",
    "卧龙电驱",
    "吉视传媒",
    "长城军工",
    "巨轮智能",
    "利欧股份",
    "阳光电源",
    "中际旭创",
    "山子高科",
    "春兴精工",
    "寒武纪-U",
    "华胜天成",
    "东方财富",
    "秦川机床",
    "国轩高科",
    "景兴纸业",
    "新易盛",
    "浙江荣泰",
    "北方稀土",
    "通富微电"
)

$colC = @(
    "华胜天成",
    "岩山科技",
    "卧龙电驱",
    "中际旭创",
    "万通发展",
    "长城军工",
    "利欧股份",
    "阳光电源",
    "三维通信",
    "北方稀土",
    "吉视传媒",
    "新易盛",
    "东方财富",
    "工业富联",
    "通富微电",
    "领益智造",
    "山子高科",
    "银之杰",
    "国轩高科",
    "巨轮智能"
)

for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
    $ws.Cells.Item($row, 3).Value = $colC[$i]
}
